$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.696.72"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "3.334.54"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.40"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.90"
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("D9").Value = "3.331.35"
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("E10").Value = "  +8.21%  "
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.99"
$ws.Range("E12").Value = "  +6.02%  "
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "691.58"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "3.873.95"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "67.710.29"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "3.327.47"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.62"
$ws.Range("E20").Value = "  +3.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("E21").Value = "  +5.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.894"
$ws.Range("E22").Value = "  +2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.49"
$ws.Range("E23").Value = "  +5.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.87"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.94"
$ws.Range("E25").Value = "  +5.42%  "
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.45"
$ws.Range("E28").Value = "  +6.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.08"
$ws.Range("E29").Value = "  +2.93%  "
$ws.Range("E30").Value = "  +4.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  +7.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "571.98"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.01"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.29"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.705.96"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.41"
$ws.Range("E39").Value = "  +14.24%  "
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("E41").Value = "  +7.95%  "
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("E43").Value = "  +3.49%  "
$ws.Range("E44").Value = "  +3.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.66"
$ws.Range("E47").Value = "  +6.99%  "
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.76"
